# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.305.43'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.590.82'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.93'
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0610'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.34'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.815.01'
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.642.83'
$ws.Range('E13').Value = '  +3.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.319.46'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  +3.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '212.00'
$ws.Range('E20').Value = '  +2.84%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.01'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('E24').Value = '  -2.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.21'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.05'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.19'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.338.25'
$ws.Range('E34').Value = '  +4.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.45'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.602'
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('E39').Value = '  -15.14%  '
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('E41').Value = '  +3.76%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.90'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.726.93'
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.93'
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0982'
$ws.Range('E49').Value = '  -2.74%  '
$ws.Range('E50').Value = '  -0.85%  '
